$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 previously held a recognition record (A4=3, B4=333, C4=date 45652).
# Clear that row's data back out (C4 keeps its date number-format style),
# as if the day's card hadn't been scanned yet.
$ws.Range("A4:C4").ClearContents()

# Reflect the new active selection (was C2, now C4).
$ws.Range("C4").Select()
